$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3.073632338554148
$ws.Cells.Item(2, 3).Value = 0.3703361064904982
$ws.Cells.Item(2, 4).Value = 0.01890383840275156
$ws.Cells.Item(2, 6).Value = 3.45912806744245
$ws.Cells.Item(2, 7).Value = 0.002586810006626408
$ws.Cells.Item(2, 9).Value = 1.960939863284167
$ws.Cells.Item(2, 10).Value = 0.1245577869196595
$ws.Cells.Item(2, 12).Value = 0.5043837092371888
$ws.Cells.Item(2, 14).Value = 2.077668157563764
$ws.Cells.Item(3, 2).Value = 2.958093105352759
$ws.Cells.Item(3, 3).Value = 0.3388644192335448
$ws.Cells.Item(3, 4).Value = 0.01848680911183287
$ws.Cells.Item(3, 6).Value = 3.444573819065312
$ws.Cells.Item(3, 7).Value = 0.002592436514923426
$ws.Cells.Item(3, 9).Value = 1.961259606711636
$ws.Cells.Item(3, 10).Value = 0.1252492598013513
$ws.Cells.Item(3, 12).Value = 0.4977875773932254
$ws.Cells.Item(3, 14).Value = 2.100866397189037
$ws.Cells.Item(4, 2).Value = 2.888837321554661
$ws.Cells.Item(4, 3).Value = 0.3197376202043927
$ws.Cells.Item(4, 4).Value = 0.01822512993758885
$ws.Cells.Item(4, 6).Value = 3.437547208888589
$ws.Cells.Item(4, 7).Value = 0.002596073570168111
$ws.Cells.Item(4, 9).Value = 1.962495882925005
$ws.Cells.Item(4, 10).Value = 0.1257000761031009
$ws.Cells.Item(4, 12).Value = 0.493981482341681
$ws.Cells.Item(4, 14).Value = 2.115833324808534
$ws.Cells.Item(5, 2).Value = 2.86103817259783
$ws.Cells.Item(5, 3).Value = 0.3119922203108558
$ws.Cells.Item(5, 4).Value = 0.01811704232074618
$ws.Cells.Item(5, 6).Value = 3.435162770495523
$ws.Cells.Item(5, 7).Value = 0.002597601712332191
$ws.Cells.Item(5, 9).Value = 1.963260515479668
$ws.Cells.Item(5, 10).Value = 0.1258903982789104
$ws.Cells.Item(5, 12).Value = 0.4924917980478085
$ws.Cells.Item(5, 14).Value = 2.122114019259399
$ws.Cells.Item(6, 2).Value = 2.856447680824374
$ws.Cells.Item(6, 3).Value = 0.3107090402963024
$ws.Cells.Item(6, 4).Value = 0.01809900570787981
$ws.Cells.Item(6, 6).Value = 3.43479572796511
$ws.Cells.Item(6, 7).Value = 0.002597858242773705
$ws.Cells.Item(6, 9).Value = 1.963403216970079
$ws.Cells.Item(6, 10).Value = 0.1259224007158455
$ws.Cells.Item(6, 12).Value = 0.4922481409673622
$ws.Cells.Item(6, 14).Value = 2.123167880765642
$ws.Cells.Item(7, 2).Value = 2.888460700207133
$ws.Cells.Item(7, 3).Value = 0.3196329657124011
$ws.Cells.Item(7, 4).Value = 0.01822367815569947
$ws.Cells.Item(7, 6).Value = 3.437513113908864
$ws.Cells.Item(7, 7).Value = 0.002596093992636226
$ws.Cells.Item(7, 9).Value = 1.96250513979173
$ws.Cells.Item(7, 10).Value = 0.1257026160740704
$ws.Cells.Item(7, 12).Value = 0.4939611436299884
$ws.Cells.Item(7, 14).Value = 2.115917293755771
$ws.Cells.Item(8, 2).Value = 3.033444019841966
$ws.Cells.Item(8, 3).Value = 0.359443421372049
$ws.Cells.Item(8, 4).Value = 0.01876119063580539
$ws.Cells.Item(8, 6).Value = 3.453712614405362
$ws.Cells.Item(8, 7).Value = 0.002588712273780003
$ws.Cells.Item(8, 9).Value = 1.960833861341015
$ws.Cells.Item(8, 10).Value = 0.124790764980852
$ws.Cells.Item(8, 12).Value = 0.5020587167572899
$ws.Cells.Item(8, 14).Value = 2.085516659256939
$ws.Cells.Item(9, 2).Value = 3.331185931307061
$ws.Cells.Item(9, 3).Value = 0.4391056434031384
$ws.Cells.Item(9, 4).Value = 0.01977237502270768
$ws.Cells.Item(9, 6).Value = 3.500696507374158
$ws.Cells.Item(9, 7).Value = 0.002575676458970181
$ws.Cells.Item(9, 9).Value = 1.965840612816891
$ws.Cells.Item(9, 10).Value = 0.1232104304758552
$ws.Cells.Item(9, 12).Value = 0.5198759791533405
$ws.Cells.Item(9, 14).Value = 2.031650522454207
$ws.Cells.Item(10, 2).Value = 3.558228390742954
$ws.Cells.Item(10, 3).Value = 0.4986550978842956
$ws.Cells.Item(10, 4).Value = 0.02049158905366966
$ws.Cells.Item(10, 6).Value = 3.544590010141093
$ws.Cells.Item(10, 7).Value = 0.002566966641655455
$ws.Cells.Item(10, 9).Value = 1.9746181520226
$ws.Cells.Item(10, 10).Value = 0.1221753764974665
$ws.Cells.Item(10, 12).Value = 0.5341532004237308
$ws.Cells.Item(10, 14).Value = 1.995595086843537
$ws.Cells.Item(11, 2).Value = 3.663341173848494
$ws.Cells.Item(11, 3).Value = 0.5259787875648385
$ws.Cells.Item(11, 4).Value = 0.02081416415385284
$ws.Cells.Item(11, 6).Value = 3.566615979968333
$ws.Cells.Item(11, 7).Value = 0.002563190555188475
$ws.Cells.Item(11, 9).Value = 1.979729463269123
$ws.Cells.Item(11, 10).Value = 0.1217317269264111
$ws.Cells.Item(11, 12).Value = 0.5409073805986253
$ws.Cells.Item(11, 14).Value = 1.979960181839239
$ws.Cells.Item(12, 2).Value = 3.703409310137488
$ws.Cells.Item(12, 3).Value = 0.5363600521629905
$ws.Cells.Item(12, 4).Value = 0.02093569682340579
$ws.Cells.Item(12, 6).Value = 3.575254335010754
$ws.Cells.Item(12, 7).Value = 0.002561787238779547
$ws.Cells.Item(12, 9).Value = 1.981826657385469
$ws.Cells.Item(12, 10).Value = 0.1215676297474833
$ws.Cells.Item(12, 12).Value = 0.543502389313403
$ws.Cells.Item(12, 14).Value = 1.974150263844777
$ws.Cells.Item(13, 2).Value = 3.694768152185759
$ws.Cells.Item(13, 3).Value = 0.5341227212744002
$ws.Cells.Item(13, 4).Value = 0.02090954947095014
$ws.Cells.Item(13, 6).Value = 3.573380649539303
$ws.Cells.Item(13, 7).Value = 0.002562088287409951
$ws.Cells.Item(13, 9).Value = 1.981367786212815
$ws.Cells.Item(13, 10).Value = 0.1216027975422023
$ws.Cells.Item(13, 12).Value = 0.5429418459536208
$ws.Cells.Item(13, 14).Value = 1.975396604548635
$ws.Cells.Item(14, 2).Value = 3.666632302259416
$ws.Cells.Item(14, 3).Value = 0.5268321675737866
$ws.Cells.Item(14, 4).Value = 0.02082417489413224
$ws.Cells.Item(14, 6).Value = 3.567320689203427
$ws.Cells.Item(14, 7).Value = 0.002563074571331914
$ws.Cells.Item(14, 9).Value = 1.979898756140727
$ws.Cells.Item(14, 10).Value = 0.1217181483718548
$ws.Cells.Item(14, 12).Value = 0.5411201247554374
$ws.Cells.Item(14, 14).Value = 1.979479975836064
$ws.Cells.Item(15, 2).Value = 3.649432721242135
$ws.Cells.Item(15, 3).Value = 0.5223709900527638
$ws.Cells.Item(15, 4).Value = 0.02077180107144905
$ws.Cells.Item(15, 6).Value = 3.56364759182236
$ws.Cells.Item(15, 7).Value = 0.002563682159067787
$ws.Cells.Item(15, 9).Value = 1.979020009526252
$ws.Cells.Item(15, 10).Value = 0.1217893121650757
$ws.Cells.Item(15, 12).Value = 0.5400091330667181
$ws.Cells.Item(15, 14).Value = 1.981995584715008
$ws.Cells.Item(16, 2).Value = 3.551395936781603
$ws.Cells.Item(16, 3).Value = 0.4968742037954712
$ws.Cells.Item(16, 4).Value = 0.02047041903145264
$ws.Cells.Item(16, 6).Value = 3.543192119383605
$ws.Cells.Item(16, 7).Value = 0.002567217149086649
$ws.Cells.Item(16, 9).Value = 1.974306690971957
$ws.Cells.Item(16, 10).Value = 0.1222049166967789
$ws.Cells.Item(16, 12).Value = 0.5337170224869254
$ws.Cells.Item(16, 14).Value = 1.996632330753393
$ws.Cells.Item(17, 2).Value = 3.491722992193274
$ws.Cells.Item(17, 3).Value = 0.4812932511466101
$ws.Cells.Item(17, 4).Value = 0.02028438259057452
$ws.Cells.Item(17, 6).Value = 3.531171716754557
$ws.Cells.Item(17, 7).Value = 0.002569433297844114
$ws.Cells.Item(17, 9).Value = 1.97170225029538
$ws.Cells.Item(17, 10).Value = 0.1224668376859768
$ws.Cells.Item(17, 12).Value = 0.5299234864269238
$ws.Cells.Item(17, 14).Value = 2.005808204079081
$ws.Cells.Item(18, 2).Value = 3.457572878152519
$ws.Cells.Item(18, 3).Value = 0.4723535156290382
$ws.Cells.Item(18, 4).Value = 0.02017694394210778
$ws.Cells.Item(18, 6).Value = 3.524451567987455
$ws.Cells.Item(18, 7).Value = 0.002570725488609922
$ws.Cells.Item(18, 9).Value = 1.970309460114706
$ws.Cells.Item(18, 10).Value = 0.1226200485505506
$ws.Cells.Item(18, 12).Value = 0.5277659578456166
$ws.Cells.Item(18, 14).Value = 2.011158082775609
$ws.Cells.Item(19, 2).Value = 3.446039778058037
$ws.Cells.Item(19, 3).Value = 0.4693304356524663
$ws.Cells.Item(19, 4).Value = 0.02014049101540394
$ws.Cells.Item(19, 6).Value = 3.522209456352869
$ws.Cells.Item(19, 7).Value = 0.00257116601624155
$ws.Cells.Item(19, 9).Value = 1.969855930963149
$ws.Cells.Item(19, 10).Value = 0.1226723631971183
$ws.Cells.Item(19, 12).Value = 0.5270396473714243
$ws.Cells.Item(19, 14).Value = 2.012981841014991
$ws.Cells.Item(20, 2).Value = 3.498057451020486
$ws.Cells.Item(20, 3).Value = 0.4829495866722482
$ws.Cells.Item(20, 4).Value = 0.02030423124861969
$ws.Cells.Item(20, 6).Value = 3.532431253828065
$ws.Cells.Item(20, 7).Value = 0.002569195572537364
$ws.Cells.Item(20, 9).Value = 1.971968602097633
$ws.Cells.Item(20, 10).Value = 0.1224386907732695
$ws.Cells.Item(20, 12).Value = 0.5303247878436821
$ws.Cells.Item(20, 14).Value = 2.004823945623571
$ws.Cells.Item(21, 2).Value = 3.674889304310398
$ws.Cells.Item(21, 3).Value = 0.5289726414971483
$ws.Cells.Item(21, 4).Value = 0.02084926795546949
$ws.Cells.Item(21, 6).Value = 3.569092557023879
$ws.Cells.Item(21, 7).Value = 0.002562784154946245
$ws.Cells.Item(21, 9).Value = 1.980325852006857
$ws.Cells.Item(21, 10).Value = 0.1216841611683694
$ws.Cells.Item(21, 12).Value = 0.5416541944507998
$ws.Cells.Item(21, 14).Value = 1.978277583595535
$ws.Cells.Item(22, 2).Value = 3.791999661789475
$ws.Cells.Item(22, 3).Value = 0.5592520627670297
$ws.Cells.Item(22, 4).Value = 0.02120188968329195
$ws.Cells.Item(22, 6).Value = 3.594788219805594
$ws.Cells.Item(22, 7).Value = 0.002558748932523246
$ws.Cells.Item(22, 9).Value = 1.986730422557542
$ws.Cells.Item(22, 10).Value = 0.1212137800318018
$ws.Cells.Item(22, 12).Value = 0.5492763292830602
$ws.Cells.Item(22, 14).Value = 1.961573329270614
$ws.Cells.Item(23, 2).Value = 3.729354343073112
$ws.Cells.Item(23, 3).Value = 0.5430727751992208
$ws.Cells.Item(23, 4).Value = 0.02101400348517224
$ws.Cells.Item(23, 6).Value = 3.580914645627018
$ws.Cells.Item(23, 7).Value = 0.002560888472099108
$ws.Cells.Item(23, 9).Value = 1.983225653933246
$ws.Cells.Item(23, 10).Value = 0.1214627525971448
$ws.Cells.Item(23, 12).Value = 0.545188315536123
$ws.Cells.Item(23, 14).Value = 1.970429521267704
$ws.Cells.Item(24, 2).Value = 3.495193154312801
$ws.Cells.Item(24, 3).Value = 0.4822007014853966
$ws.Cells.Item(24, 4).Value = 0.02029525918032427
$ws.Cells.Item(24, 6).Value = 3.531861223743675
$ws.Cells.Item(24, 7).Value = 0.002569302991709379
$ws.Cells.Item(24, 9).Value = 1.971847858899793
$ws.Cells.Item(24, 10).Value = 0.1224514078030632
$ws.Cells.Item(24, 12).Value = 0.5301432865100679
$ws.Cells.Item(24, 14).Value = 2.00526869670469
$ws.Cells.Item(25, 2).Value = 3.24919047076753
$ws.Cells.Item(25, 3).Value = 0.4173791349065255
$ws.Cells.Item(25, 4).Value = 0.0195031630458562
$ws.Cells.Item(25, 6).Value = 3.486347040646393
$ws.Cells.Item(25, 7).Value = 0.002579049906938861
$ws.Cells.Item(25, 9).Value = 1.963594279161981
$ws.Cells.Item(25, 10).Value = 0.1236157752277496
$ws.Cells.Item(25, 12).Value = 0.5148479177941994
$ws.Cells.Item(25, 14).Value = 2.045605704023902
